$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove column F (and its header) entirely -- data now ends at column E
$ws.Columns("F").Delete()

# Bulk-write A2:E61 with the refreshed monthly data (rotated month order +
# two additional yearly cycles appended: 2022 and 2023)
$data = New-Object "object[,]" 60,5
$data[0,0] = "2018-10"
$data[0,1] = 0
$data[0,2] = 0
$data[0,3] = 0
$data[0,4] = 0
$data[1,0] = "2018-11"
$data[1,1] = 0
$data[1,2] = 0
$data[1,3] = 0
$data[1,4] = 0
$data[2,0] = "2018-12"
$data[2,1] = 0
$data[2,2] = 0
$data[2,3] = 0
$data[2,4] = 0
$data[3,0] = "2018-02"
$data[3,1] = "'"
$data[3,2] = "'"
$data[3,3] = 0
$data[3,4] = 0
$data[4,0] = "2018-03"
$data[4,1] = 0
$data[4,2] = 0
$data[4,3] = 0
$data[4,4] = 0
$data[5,0] = "2018-04"
$data[5,1] = 0
$data[5,2] = 0
$data[5,3] = 0
$data[5,4] = 0
$data[6,0] = "2018-05"
$data[6,1] = 0
$data[6,2] = 0
$data[6,3] = 0
$data[6,4] = 0
$data[7,0] = "2018-06"
$data[7,1] = 0
$data[7,2] = 0
$data[7,3] = 0
$data[7,4] = 0
$data[8,0] = "2018-07"
$data[8,1] = 0
$data[8,2] = 0
$data[8,3] = 0
$data[8,4] = 0
$data[9,0] = "2018-08"
$data[9,1] = 0
$data[9,2] = 0
$data[9,3] = 0
$data[9,4] = 0
$data[10,0] = "2018-09"
$data[10,1] = 31.5
$data[10,2] = 0.6
$data[10,3] = 3.9
$data[10,4] = 19.8
$data[11,0] = "2019-10"
$data[11,1] = 0
$data[11,2] = 0
$data[11,3] = 0
$data[11,4] = 0
$data[12,0] = "2019-11"
$data[12,1] = 0
$data[12,2] = 0
$data[12,3] = 0
$data[12,4] = 0
$data[13,0] = "2019-12"
$data[13,1] = 0
$data[13,2] = 0
$data[13,3] = 0
$data[13,4] = 0
$data[14,0] = "2019-02"
$data[14,1] = "'"
$data[14,2] = "'"
$data[14,3] = 0
$data[14,4] = 0
$data[15,0] = "2019-03"
$data[15,1] = 0
$data[15,2] = 0
$data[15,3] = 0
$data[15,4] = 0
$data[16,0] = "2019-04"
$data[16,1] = 0
$data[16,2] = 0
$data[16,3] = 0
$data[16,4] = 0
$data[17,0] = "2019-05"
$data[17,1] = 0
$data[17,2] = 0
$data[17,3] = 0
$data[17,4] = 0
$data[18,0] = "2019-06"
$data[18,1] = 0
$data[18,2] = 0
$data[18,3] = 0
$data[18,4] = 0
$data[19,0] = "2019-07"
$data[19,1] = 0
$data[19,2] = 0
$data[19,3] = 0
$data[19,4] = 0
$data[20,0] = "2019-08"
$data[20,1] = 0
$data[20,2] = 0
$data[20,3] = 0
$data[20,4] = 0
$data[21,0] = "2019-09"
$data[21,1] = 0
$data[21,2] = 0
$data[21,3] = 0
$data[21,4] = 0
$data[22,0] = "2020-10"
$data[22,1] = 0
$data[22,2] = 0
$data[22,3] = 0
$data[22,4] = 1777.3
$data[23,0] = "2020-11"
$data[23,1] = 0
$data[23,2] = 0
$data[23,3] = 0
$data[23,4] = 1777.3
$data[24,0] = "2020-12"
$data[24,1] = 0
$data[24,2] = 0
$data[24,3] = 0
$data[24,4] = -100
$data[25,0] = "2020-02"
$data[25,1] = "'"
$data[25,2] = "'"
$data[25,3] = 0
$data[25,4] = 0
$data[26,0] = "2020-03"
$data[26,1] = 0
$data[26,2] = 0
$data[26,3] = 0
$data[26,4] = 0
$data[27,0] = "2020-04"
$data[27,1] = 0
$data[27,2] = 0
$data[27,3] = 0
$data[27,4] = 0
$data[28,0] = "2020-05"
$data[28,1] = 0
$data[28,2] = 0
$data[28,3] = 0
$data[28,4] = 0
$data[29,0] = "2020-06"
$data[29,1] = -100
$data[29,2] = 0
$data[29,3] = 0
$data[29,4] = 1777.3
$data[30,0] = "2020-07"
$data[30,1] = 462.7
$data[30,2] = 0.1
$data[30,3] = 0.2
$data[30,4] = -30.4
$data[31,0] = "2020-08"
$data[31,1] = 233.3
$data[31,2] = 0.1
$data[31,3] = 0.2
$data[31,4] = -9.1
$data[32,0] = "2020-09"
$data[32,1] = 0
$data[32,2] = 0
$data[32,3] = 0
$data[32,4] = 1777.3
$data[33,0] = "2021-10"
$data[33,1] = 0
$data[33,2] = 0
$data[33,3] = 0
$data[33,4] = 0
$data[34,0] = "2021-11"
$data[34,1] = 0
$data[34,2] = 0
$data[34,3] = 0
$data[34,4] = 0
$data[35,0] = "2021-12"
$data[35,1] = 0
$data[35,2] = 0
$data[35,3] = 0
$data[35,4] = 0
$data[36,0] = "2021-02"
$data[36,1] = "'"
$data[36,2] = "'"
$data[36,3] = 0
$data[36,4] = -100
$data[37,0] = "2021-03"
$data[37,1] = 0
$data[37,2] = 0
$data[37,3] = 0
$data[37,4] = -100
$data[38,0] = "2021-04"
$data[38,1] = -100
$data[38,2] = 0
$data[38,3] = 0
$data[38,4] = -100
$data[39,0] = "2021-05"
$data[39,1] = 0
$data[39,2] = 0
$data[39,3] = 0
$data[39,4] = -100
$data[40,0] = "2021-06"
$data[40,1] = 0
$data[40,2] = 0
$data[40,3] = 0
$data[40,4] = -100
$data[41,0] = "2021-07"
$data[41,1] = 0
$data[41,2] = 0
$data[41,3] = 0
$data[41,4] = 0
$data[42,0] = "2021-08"
$data[42,1] = 0
$data[42,2] = 0
$data[42,3] = 0
$data[42,4] = 0
$data[43,0] = "2021-09"
$data[43,1] = 0
$data[43,2] = 0
$data[43,3] = 0
$data[43,4] = 0
$data[44,0] = "2022-10"
$data[44,1] = 0
$data[44,2] = 0
$data[44,3] = 0
$data[44,4] = 0
$data[45,0] = "2022-11"
$data[45,1] = 0
$data[45,2] = 0
$data[45,3] = 0
$data[45,4] = 0
$data[46,0] = "2022-12"
$data[46,1] = 0
$data[46,2] = 0
$data[46,3] = 0
$data[46,4] = 0
$data[47,0] = "2022-02"
$data[47,1] = "'"
$data[47,2] = "'"
$data[47,3] = 0
$data[47,4] = 0
$data[48,0] = "2022-03"
$data[48,1] = 0
$data[48,2] = 0
$data[48,3] = 0
$data[48,4] = 0
$data[49,0] = "2022-04"
$data[49,1] = 0
$data[49,2] = 0
$data[49,3] = 0
$data[49,4] = 0
$data[50,0] = "2022-05"
$data[50,1] = 0
$data[50,2] = 0
$data[50,3] = 0
$data[50,4] = 0
$data[51,0] = "2022-06"
$data[51,1] = 0
$data[51,2] = 0
$data[51,3] = 0
$data[51,4] = 0
$data[52,0] = "2022-07"
$data[52,1] = 0
$data[52,2] = 0
$data[52,3] = 0
$data[52,4] = 0
$data[53,0] = "2022-08"
$data[53,1] = 0
$data[53,2] = 0
$data[53,3] = 0
$data[53,4] = 0
$data[54,0] = "2022-09"
$data[54,1] = 0
$data[54,2] = 0
$data[54,3] = 0
$data[54,4] = 0
$data[55,0] = "2023-02"
$data[55,1] = "'"
$data[55,2] = "'"
$data[55,3] = 0
$data[55,4] = 0
$data[56,0] = "2023-03"
$data[56,1] = 0
$data[56,2] = 0
$data[56,3] = 0
$data[56,4] = 0
$data[57,0] = "2023-04"
$data[57,1] = 0
$data[57,2] = 0
$data[57,3] = 0
$data[57,4] = 0
$data[58,0] = "2023-05"
$data[58,1] = 0
$data[58,2] = 0
$data[58,3] = 0
$data[58,4] = 0
$data[59,0] = "2023-06"
$data[59,1] = 0
$data[59,2] = 0
$data[59,3] = 0
$data[59,4] = 0
$ws.Range("A2:E61").Value2 = $data

# Carry the header-column style (bold, bordered, centered) onto the newly
# added rows 46:61, matching the style already used by rows 2:45
$ws.Range("A2").Copy()
$ws.Range("A46:A61").PasteSpecial(-4122)
$excel.CutCopyMode = $false
